$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-5 (columns B:F) with the new simulated data
$ws.Cells.Item(2, 2).Value = 0.4524751023004333
$ws.Cells.Item(2, 3).Value = 0.7399800961295603
$ws.Cells.Item(2, 4).Value = 0.959507423691378
$ws.Cells.Item(2, 5).Value = 0.9795444980660031
$ws.Cells.Item(2, 6).Value = 0.901572637837064

$ws.Cells.Item(3, 2).Value = 0.06200830852274453
$ws.Cells.Item(3, 3).Value = 0.6054022942997469
$ws.Cells.Item(3, 4).Value = 0.6570257936368813
$ws.Cells.Item(3, 5).Value = 0.810571276099074
$ws.Cells.Item(3, 6).Value = 0.8519133780372987
$ws.Cells.Item(3, 7).Value = 10

$ws.Cells.Item(4, 2).Value = -0.06427704427340604
$ws.Cells.Item(4, 3).Value = 0.2800168996375158
$ws.Cells.Item(4, 4).Value = 0.1697350109970526
$ws.Cells.Item(4, 5).Value = 0.411989090871412
$ws.Cells.Item(4, 6).Value = 0.445784888810548
$ws.Cells.Item(4, 7).Value = 6

$ws.Cells.Item(5, 2).Value = 0.2218614552644835
$ws.Cells.Item(5, 3).Value = 0.2564559007953193
$ws.Cells.Item(5, 4).Value = 0.1149921343848131
$ws.Cells.Item(5, 5).Value = 0.3391049017410587
$ws.Cells.Item(5, 6).Value = 0.3626834130553496
$ws.Cells.Item(5, 7).Value = 2

# Delete rows 6 through 9 (Q4-Q7 data no longer needed)
$ws.Range("A6:G9").Delete()
